# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 390
    5  = 11537
    6  = 771
    7  = 111
    11 = 166
    14 = 49
    17 = 327
    18 = 1312
    19 = 70
    20 = 897
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
